$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 518; existing rows 518:591 shift down to 519:592
$ws.Rows.Item(518).Insert()

# Populate the newly inserted row 518 with its values
$ws.Cells.Item(518, 1).Value = 3
$ws.Cells.Item(518, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(518, 3).Value = "Coquimbo"
$ws.Cells.Item(518, 4).Value = "2023-07-20"
$ws.Cells.Item(518, 5).Value = 5
$ws.Cells.Item(518, 6).Value = 100114013
$ws.Cells.Item(518, 7).Value = "Zanahoria"
$ws.Cells.Item(518, 8).Value = "Sin especificar"
$ws.Cells.Item(518, 9).Value = "Primera"
$ws.Cells.Item(518, 10).Value = 230
$ws.Cells.Item(518, 11).Value = 7000
$ws.Cells.Item(518, 12).Value = 7500
$ws.Cells.Item(518, 13).Value = 7261
$ws.Cells.Item(518, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(518, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(518, 16).Value = 363
$ws.Cells.Item(518, 17).Value = 20
$ws.Cells.Item(518, 18).Value = "Hortaliza"
